$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D5").Value = "2016-01-18 12:52:19"
$wsZhCn.Range("G5").Value = "2016-01-18 12:53:01"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D5").Value = "2016-01-18 12:52:28"
$wsDeDe.Range("G5").Value = "2016-01-18 12:53:18"
